$wb = $excel.ActiveWorkbook

# --- LaLiga: append 19 new match rows starting at row 271 ---
$ws = $wb.Worksheets.Item('LaLiga')
$LaLigaData = @(
    ,('Getafe','Valencia',0,0)
    ,('Villarreal','Celta Vigo',1,0)
    ,('Granada','Elche',0,1)
    ,('Levante','Espanyol',1,1)
    ,('Barcelona','Osasuna',4,0)
    ,('Real Sociedad','Alaves',1,0)
    ,('Betis','Ath. Bilbao',1,0)
    ,('Vallecano','Sevilla',1,1)
    ,('Mallorca','Real Madryt',0,3)
    ,('Ath. Bilbao','Getafe',1,1)
    ,('Vallecano','Atl. Madryt',0,1)
    ,('Osasuna','Levante',3,1)
    ,('Elche','Valencia',0,1)
    ,('Alaves','Granada',2,3)
    ,('Real Madryt','Barcelona',0,4)
    ,('Sevilla','Real Sociedad',0,0)
    ,('Celta Vigo','Betis',0,0)
    ,('Cadiz','Villarreal',1,0)
    ,('Espanyol','Mallorca',1,0)
)
$startRow = 271
for ($i = 0; $i -lt $LaLigaData.Count; $i++) {
    $r = $startRow + $i
    $match = $LaLigaData[$i]
    $ws.Cells.Item($r, 1).Value = $match[0]
    $ws.Cells.Item($r, 2).Value = $match[1]
    $ws.Cells.Item($r, 3).Value = $match[2]
    $ws.Cells.Item($r, 4).Value = $match[3]
}

# --- PremierLeague: append 17 new match rows starting at row 273 ---
$ws = $wb.Worksheets.Item('PremierLeague')
$PremierLeagueData = @(
    ,('Manchester Utd','Tottenham',3,2)
    ,('Brentford','Burnley',2,0)
    ,('Brighton','Liverpool',0,2)
    ,('Arsenal','Leicester',2,0)
    ,('West Ham','Aston Villa',2,1)
    ,('Southampton','Watford',1,2)
    ,('Leeds','Norwich',2,1)
    ,('Everton','Wolves',0,1)
    ,('Chelsea','Newcastle',1,0)
    ,('Crystal Palace','Manchester City',0,0)
    ,('Arsenal','Liverpool',0,2)
    ,('Brighton','Tottenham',0,2)
    ,('Everton','Newcastle',1,0)
    ,('Wolves','Leeds',2,3)
    ,('Aston Villa','Arsenal',0,1)
    ,('Tottenham','West Ham',3,1)
    ,('Leicester','Brentford',2,1)
)
$startRow = 273
for ($i = 0; $i -lt $PremierLeagueData.Count; $i++) {
    $r = $startRow + $i
    $match = $PremierLeagueData[$i]
    $ws.Cells.Item($r, 1).Value = $match[0]
    $ws.Cells.Item($r, 2).Value = $match[1]
    $ws.Cells.Item($r, 3).Value = $match[2]
    $ws.Cells.Item($r, 4).Value = $match[3]
}

# --- Championship: append 20 new match rows starting at row 425 ---
$ws = $wb.Worksheets.Item('Championship')
$ChampionshipData = @(
    ,('Bournemouth','Reading',1,1)
    ,('Blackburn','Derby',3,1)
    ,('Birmingham','Middlesbrough',0,2)
    ,('Barnsley','Bristol City',2,0)
    ,('Peterborough','Swansea',2,3)
    ,('Nottingham','QPR',3,1)
    ,('Millwall','Huddersfield',2,0)
    ,('Luton','Preston',4,0)
    ,('Coventry','Hull',0,2)
    ,('Cardiff','Stoke',2,1)
    ,('Blackpool','Sheffield Utd',0,0)
    ,('Swansea','Birmingham',0,0)
    ,('Stoke','Millwall',2,0)
    ,('Reading','Blackburn',1,0)
    ,('Hull','Luton',1,3)
    ,('Huddersfield','Bournemouth',0,3)
    ,('Bristol City','West Brom',2,2)
    ,('Sheffield Utd','Barnsley',2,0)
    ,('Derby','Coventry',1,1)
    ,('QPR','Peterborough',1,3)
)
$startRow = 425
for ($i = 0; $i -lt $ChampionshipData.Count; $i++) {
    $r = $startRow + $i
    $match = $ChampionshipData[$i]
    $ws.Cells.Item($r, 1).Value = $match[0]
    $ws.Cells.Item($r, 2).Value = $match[1]
    $ws.Cells.Item($r, 3).Value = $match[2]
    $ws.Cells.Item($r, 4).Value = $match[3]
}

# --- SerieA: append 20 new match rows starting at row 277 ---
$ws = $wb.Worksheets.Item('SerieA')
$SerieAData = @(
    ,('AC Milan','Empoli',1,0)
    ,('Sampdoria','Juventus',1,3)
    ,('Spezia','Cagliari',2,0)
    ,('Salernitana','Sassuolo',2,2)
    ,('Torino','Inter',1,1)
    ,('Udinese','AS Roma',1,1)
    ,('Atalanta','Genoa',0,0)
    ,('Verona','Napoli',1,2)
    ,('Fiorentina','Bologna',1,0)
    ,('Lazio','Venezia',1,0)
    ,('Genoa','Torino',1,0)
    ,('Sassuolo','Spezia',4,1)
    ,('Cagliari','AC Milan',0,1)
    ,('Inter','Fiorentina',1,1)
    ,('Napoli','Udinese',2,1)
    ,('Bologna','Atalanta',0,1)
    ,('AS Roma','Lazio',3,0)
    ,('Juventus','Salernitana',2,0)
    ,('Empoli','Verona',1,1)
    ,('Venezia','Sampdoria',0,2)
)
$startRow = 277
for ($i = 0; $i -lt $SerieAData.Count; $i++) {
    $r = $startRow + $i
    $match = $SerieAData[$i]
    $ws.Cells.Item($r, 1).Value = $match[0]
    $ws.Cells.Item($r, 2).Value = $match[1]
    $ws.Cells.Item($r, 3).Value = $match[2]
    $ws.Cells.Item($r, 4).Value = $match[3]
}

# --- Ligue1: append 19 new match rows starting at row 273 ---
$ws = $wb.Worksheets.Item('Ligue1')
$Ligue1Data = @(
    ,('Troyes','Nantes',1,0)
    ,('Montpellier','Nice',0,0)
    ,('Brest','Marsylia',1,4)
    ,('Lyon','Rennes',2,4)
    ,('Strasbourg','Monaco',1,0)
    ,('Metz','Lens',0,0)
    ,('Clermont','Lorient',0,2)
    ,('Angers','Reims',0,1)
    ,('PSG','Bordeaux',3,0)
    ,('St. Etienne','Troyes',1,1)
    ,('Nantes','Lille',0,1)
    ,('Lens','Clermont',3,1)
    ,('Marsylia','Nice',2,1)
    ,('Reims','Lyon',0,0)
    ,('Rennes','Metz',6,1)
    ,('Lorient','Strasbourg',0,0)
    ,('Bordeaux','Montpellier',0,2)
    ,('Angers','Brest',1,0)
    ,('Monaco','PSG',3,0)
)
$startRow = 273
for ($i = 0; $i -lt $Ligue1Data.Count; $i++) {
    $r = $startRow + $i
    $match = $Ligue1Data[$i]
    $ws.Cells.Item($r, 1).Value = $match[0]
    $ws.Cells.Item($r, 2).Value = $match[1]
    $ws.Cells.Item($r, 3).Value = $match[2]
    $ws.Cells.Item($r, 4).Value = $match[3]
}
